$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray/erroneous cells (naive component forecaster bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected forecast values
$ws.Range("E3").Value = 8.262942840582976
$ws.Range("C4").Value = 6.277541464866965
$ws.Range("E4").Value = 7.915558093865016
$ws.Range("C5").Value = 6.535114773304795
$ws.Range("E7").Value = 3.624426704091555
$ws.Range("E8").Value = 3.771815305047843
$ws.Range("E10").Value = 2.56219956496937
$ws.Range("C11").Value = 2.508469427909921
$ws.Range("E11").Value = 2.632055757778873
$ws.Range("E13").Value = 1.467147844249128
$ws.Range("C14").Value = 3.047037961814514
$ws.Range("E14").Value = 2.75705424928776
$ws.Range("C15").Value = -0.22288476972816
$ws.Range("E15").Value = 1.7415595764392
$ws.Range("C16").Value = -1.165854108406639
$ws.Range("E17").Value = 2.055357398179125
$ws.Range("E18").Value = 0.9064937165318865
$ws.Range("C19").Value = 2.039329803030099
